# Adapt column header formatting to respective input-file names:
#   *_old -> *_FV2210   (left-hand comparison block)
#   *_new -> *_FV2304   (right-hand comparison block)
# Then turn the header row + data into a real Excel Table and freeze the
# header row so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldSuffix = "_old"
$newSuffix = "_new"
$fv2210    = "_FV2210"
$fv2304    = "_FV2304"

# --- 1) Rename the header row (row 1) shared-string values -----------------
$usedRange = $ws.UsedRange
$lastCol = $usedRange.Columns.Count
$lastRow = $usedRange.Rows.Count

for ($col = 1; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $val = $cell.Value()
    if ($null -eq $val) { continue }

    if ($val -like "*$oldSuffix") {
        $cell.Value = $val -replace ([regex]::Escape($oldSuffix) + '$'), $fv2210
    } elseif ($val -like "*$newSuffix") {
        $cell.Value = $val -replace ([regex]::Escape($newSuffix) + '$'), $fv2304
    }
}

# --- 2) Turn A1:<lastCol><lastRow> into an Excel Table ----------------------
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"

# --- 3) Freeze the header row (split/freeze at row 2, i.e. above row 2) ----
$ws.Range("A2").Select() | Out-Null
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select() | Out-Null
